# Move the `assignment_date` merge field from the `casa_case` data block to
# the `volunteer` data block.
#
# The field in question looks like (OOXML):
#   begin -> instrText " MERGEFIELD =casa_case.assignment_date \* MERGEFORMAT "
#         -> separate  -> cached result text "«=casa_case.assignment_date»"
#         -> end
#
# We need it to become:
#   begin -> instrText " MERGEFIELD =volunteer.assignment_date \* MERGEFORMAT "
#         -> separate  -> cached result text "«=volunteer.assignment_date»"
#         -> end

$d = $word.ActiveDocument

# --- 1. Locate the target MERGEFIELD by its field code text -----------------
$target = $null
$fields = $d.Fields
for ($i = 1; $i -le $fields.Count; $i++) {
    $candidate = $fields.Item($i)
    if ($candidate.Code.Text -match "casa_case\.assignment_date") {
        $target = $candidate
        break
    }
}

if ($target -ne $null) {
    # --- 2. Rewrite the field's instruction (the hidden `w:instrText` run) --
    # Assigning straight to `Field.Code` (rather than `Field.Code.Text`)
    # rewrites the underlying field-code runs for the field.
    $target.Code = " MERGEFIELD =volunteer.assignment_date \* MERGEFORMAT "
}

# --- 3. Rewrite the cached/displayed field result (the visible `w:t` run) --
# This text is part of the document's visible story, so plain Find & Replace
# reaches it (field instruction text above is hidden and is not reachable
# this way, which is why step 2 is handled separately through the Field
# object itself).
$null = $d.Content.Find.Execute(
    "«=casa_case.assignment_date»",  # FindText
    $true,                            # MatchCase
    $false,                           # MatchWholeWord
    $false,                           # MatchWildcards
    $false,                           # MatchSoundsLike
    $false,                           # MatchAllWordForms
    $true,                            # Forward
    1,                                # Wrap (wdFindContinue)
    $false,                           # Format
    "«=volunteer.assignment_date»",  # ReplaceWith
    2                                 # Replace (wdReplaceAll)
)

Write-Output "done"
